$wb = $excel.ActiveWorkbook

$wsSUSALUD = $wb.Worksheets.Item("SUSALUD")
$wsGestores = $wb.Worksheets.Item("Gestores")
$wsCiudadanos = $wb.Worksheets.Item("Ciudadanos")

# SUSALUD sheet: move the cell selection from B4 to D3
$wsSUSALUD.Activate()
$wsSUSALUD.Range("D3").Select()

# Gestores sheet: move the cell selection from F4 to F3
# (it stops being the tab-selected sheet once another sheet is activated later)
$wsGestores.Activate()
$wsGestores.Range("F3").Select()

# Ciudadanos sheet: fill in the two missing "reason" cells for rows 4 and 5
# using the same text already used in the other columns of those rows
$wsCiudadanos.Range("B4").Value = "Nueva Solicitud"
$wsCiudadanos.Range("B5").Value = "Estadisticas"

# Ciudadanos becomes the active/tab-selected sheet, with selection on C6
$wsCiudadanos.Activate()
$wsCiudadanos.Range("C6").Select()
